$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (B1:T1): relabel columns so "KilometersDemand" sits right after
# "Kilometersopex" (commit: "put sinkDemand in outcolumns"); everything after shifts down one.
$headers = @(
    "Total Cost",
    "crudeoil",
    "hydrogen",
    "biomass",
    "CrOilopex",
    "H2opex",
    "BMopex",
    "Kilometersopex",
    "KilometersDemand",
    "Refineryopex",
    "RefineryTotalEff",
    "MtGopex",
    "MtGTotalEff",
    "Gtkmopex",
    "GtkmTotalEff",
    "B2gasopex",
    "B2gasTotalEff",
    "GasHubopex",
    "KmHubopex"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# --- Data rows 2-11: refreshed simulation output
# Row 2
$ws.Range("B2").Value = 158.0258831310304
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 2481.118429995444
$ws.Range("E2").Value = 999.9999999999998
$ws.Range("F2").Value = 0.03220181821443927
$ws.Range("G2").Value = 0.06217450458457861
$ws.Range("H2").Value = 0.003763573930396155
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 984.1866627823251
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0.8981465570493671
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0.8377493818800267
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0.3865973522576316
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0.4672112613902339
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0

# Row 3
$ws.Range("B3").Value = 110.3849719851686
$ws.Range("C3").Value = 546.448087431694
$ws.Range("D3").Value = 2162.013968557432
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0.006873465353439534
$ws.Range("G3").Value = 0.04931928356591516
$ws.Range("H3").Value = 0.01706439378853211
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 996.1270955509859
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0.919743317273846
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0.8310796678379754
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0.4332121582601848
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0.5226948216421022
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0

# Row 4
$ws.Range("B4").Value = 124.3593793844814
$ws.Range("C4").Value = 546.448087431694
$ws.Range("D4").Value = 2653.781307140937
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0.0172744924564318
$ws.Range("G4").Value = 0.04330415837623047
$ws.Range("H4").Value = 0.03768709715816157
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 1000.006194004161
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0.9328441864123959
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0.8309250423501848
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0.3683475351869421
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0.5339359156569512
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0

# Row 5
$ws.Range("B5").Value = 141.8437277552752
$ws.Range("C5").Value = 546.448087431694
$ws.Range("D5").Value = 2356.328596929176
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0.008806884639781748
$ws.Range("G5").Value = 0.05815454714856717
$ws.Range("H5").Value = 0.04994318697412082
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 987.3263226097781
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0.9140028381967513
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0.8637779980827724
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0.389508584947223
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0.5016422672066417
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 0

# Row 6
$ws.Range("B6").Value = 145.7292557456382
$ws.Range("C6").Value = 546.4480874316939
$ws.Range("D6").Value = 2609.746007355325
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0.03010279422794272
$ws.Range("G6").Value = 0.04953725038722792
$ws.Range("H6").Value = 0.01553751704643268
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 1000.017253009715
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0.9313166937594809
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0.8695850964250099
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0.3599369284663064
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = 0.4383399087229943
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 0

# Row 7
$ws.Range("B7").Value = 127.8121957657857
$ws.Range("C7").Value = 546.4480874316939
$ws.Range("D7").Value = 2304.994074539171
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0.00120591452773015
$ws.Range("G7").Value = 0.05516423121561503
$ws.Range("H7").Value = 0.02393140519025732
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 989.6570582566266
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0.9266935162889555
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0.8374978788903947
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0.4061268596462472
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = 0.5967957834561863
$ws.Range("S7").Value = 0
$ws.Range("T7").Value = 0

# Row 8
$ws.Range("B8").Value = 156.1305081592099
$ws.Range("C8").Value = 546.448087431694
$ws.Range("D8").Value = 2740.856160060153
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0.007362810529836793
$ws.Range("G8").Value = 0.05549620466902587
$ws.Range("H8").Value = 0.01897901985493813
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 986.8091520116678
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0.941356053461611
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0.8302975069283086
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0.3536787785515837
$ws.Range("Q8").Value = 0
$ws.Range("R8").Value = 0.465719931860492
$ws.Range("S8").Value = 0
$ws.Range("T8").Value = 0

# Row 9
$ws.Range("B9").Value = 83.14618530030624
$ws.Range("C9").Value = 546.448087431694
$ws.Range("D9").Value = 2908.268621478813
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0.005837351916810223
$ws.Range("G9").Value = 0.02749277522687808
$ws.Range("H9").Value = 0.03562161052151295
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 999.7737159298745
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0.9412176756403355
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0.8556829790046292
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0.3329380011120837
$ws.Range("Q9").Value = 0
$ws.Range("R9").Value = 0.5177928757904908
$ws.Range("S9").Value = 0
$ws.Range("T9").Value = 0

# Row 10
$ws.Range("B10").Value = 131.619468934875
$ws.Range("C10").Value = 546.448087431694
$ws.Range("D10").Value = 2441.13487897873
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0.002649223250108834
$ws.Range("G10").Value = 0.05332429890606127
$ws.Range("H10").Value = 0.02307426654619004
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 990.8067002548242
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0.904286079009541
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0.842007929763322
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0.3886125566011299
$ws.Range("Q10").Value = 0
$ws.Range("R10").Value = 0.565683438343288
$ws.Range("S10").Value = 0
$ws.Range("T10").Value = 0

# Row 11
$ws.Range("B11").Value = 73.62817470455298
$ws.Range("C11").Value = 546.448087431694
$ws.Range("D11").Value = 2643.581028045432
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0.002983402884718083
$ws.Range("G11").Value = 0.02723498888074303
$ws.Range("H11").Value = 0.01989741273540744
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 992.5833045965045
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0.9143952679391637
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0.8386706953929792
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0.3653550668953141
$ws.Range("Q11").Value = 0
$ws.Range("R11").Value = 0.5470094928806483
$ws.Range("S11").Value = 0
$ws.Range("T11").Value = 0

# --- Column widths (best-fit to new content)
$ws.Columns.Item(1).ColumnWidth = 1.1666666666666665
$ws.Columns.Item(2).ColumnWidth = 11.166666666666666
$ws.Columns.Item(3).ColumnWidth = 11.166666666666666
$ws.Columns.Item(4).ColumnWidth = 11.166666666666666
$ws.Columns.Item(5).ColumnWidth = 7.451822916666666
$ws.Columns.Item(6).ColumnWidth = 11.166666666666666
$ws.Columns.Item(7).ColumnWidth = 11.166666666666666
$ws.Columns.Item(8).ColumnWidth = 11.166666666666666
$ws.Columns.Item(9).ColumnWidth = 14.451822916666666
$ws.Columns.Item(10).ColumnWidth = 17.592447916666668
$ws.Columns.Item(11).ColumnWidth = 12.307291666666666
$ws.Columns.Item(12).ColumnWidth = 14.736979166666666
$ws.Columns.Item(13).ColumnWidth = 8.592447916666666
$ws.Columns.Item(14).ColumnWidth = 11.166666666666666
$ws.Columns.Item(15).ColumnWidth = 9.451822916666666
$ws.Columns.Item(16).ColumnWidth = 11.877604166666666
$ws.Columns.Item(17).ColumnWidth = 9.592447916666666
$ws.Columns.Item(18).ColumnWidth = 12.022135416666666
$ws.Columns.Item(19).ColumnWidth = 11.451822916666666
$ws.Columns.Item(20).ColumnWidth = 11.022135416666666

# --- Restore active selection to I13
$ws.Range("I13").Select()
